$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Helper: locate the paragraph whose text contains a given marker substring.
# ---------------------------------------------------------------------------
function Find-ParagraphIndex($marker) {
    $count = $d.Paragraphs.Count
    for ($i = 1; $i -le $count; $i++) {
        $p = $d.Paragraphs($i)
        if ($p.Range.Text.Contains($marker)) {
            return $i
        }
    }
    return -1
}

# ===========================================================================
# Change 1: split the "singleton pattern was dropped" paragraph and add a
# new list item ("Added list of ignored phrases ...") right after it,
# carrying the _GoBack bookmark into the middle of the new paragraph's text.
# ===========================================================================
$idx1 = Find-ParagraphIndex("singleton pattern was dropped")
$p1 = $d.Paragraphs($idx1)
$r1 = $p1.Range
$fullText = $r1.Text
$textLen = $fullText.Length
# Range.Text includes the trailing paragraph mark (chr 13) as its last
# character, so the "real" end of the visible text is one character short
# of the range end.
$splitPos = $r1.Start + ($textLen - 1)

# Remove the existing _GoBack bookmark - it will be re-created further down,
# inside the text of the newly inserted paragraph.
$d.Bookmarks("_GoBack").Delete()

# Split the paragraph in two; the new (second) paragraph inherits the same
# paragraph formatting (style / numbering / indent) automatically.
$splitRange = $d.Range($splitPos, $splitPos)
$null = $splitRange.InsertParagraphAfter()

$newParaIndex = $idx1 + 1
$newPara = $d.Paragraphs($newParaIndex)
$newParaRange = $newPara.Range
$insertStart = $newParaRange.Start

$quoteOpen = [char]0x201C
$quoteClose = [char]0x201D
$textBeforeBookmark = "Added list of ignored phrases to the Interaction classes, to be able to, for example,"
$textAfterBookmark = " type " + $quoteOpen + "go to" + $quoteClose + " instead of " + $quoteOpen + "go" + $quoteClose + ", and achieve the same result"

$insRange = $d.Range($insertStart, $insertStart)
$null = $insRange.InsertAfter($textBeforeBookmark + $textAfterBookmark)

# Re-anchor the bookmark exactly between the two pieces of text using a
# Range derived (via Duplicate/Move) from the live paragraph Range - building
# a brand-new Range object at a raw offset right after an insertion is not
# reliable for Bookmarks.Add in this host.
$newPara2 = $d.Paragraphs($newParaIndex)
$bmRange = $newPara2.Range.Duplicate()
$null = $bmRange.MoveStart(1, $textBeforeBookmark.Length)
$bmRange.Collapse(1)
$null = $d.Bookmarks.Add("_GoBack", $bmRange)

# ===========================================================================
# Change 2: move the <w:lastRenderedPageBreak/> marker from the
# "Create mapping annotations ..." paragraph to the preceding
# "two locations adjacent ..." paragraph. This element sits directly before
# the <w:t> inside the paragraph's (only) run, so the safest way to relocate
# it precisely is to rewrite each paragraph's WordprocessingML in place via
# Range.InsertXML (which replaces exactly the content of the Range it is
# called on).
# ===========================================================================
function Set-ParagraphOpenXml($paraIndex, $innerBodyXml) {
    $para = $d.Paragraphs($paraIndex)
    $range = $para.Range
    $pkg = '<?xml version="1.0"?>' +
        '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
        '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
        '<pkg:xmlData>' +
        '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
        '<w:body>' + $innerBodyXml + '</w:body>' +
        '</w:document>' +
        '</pkg:xmlData></pkg:part></pkg:package>'
    $null = $range.InsertXML($pkg)
}

$idxTwoLocations = Find-ParagraphIndex("two locations adjacent to each other")
$twoLocationsXml = '<w:p w:rsidR="008B0220" w:rsidRPr="00107C01" w:rsidRDefault="008B0220" w:rsidP="00107C01">' +
    '<w:pPr><w:pStyle w:val="Listenabsatz"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="21"/></w:numPr>' +
    '<w:ind w:left="1068"/><w:rPr><w:lang w:val="en-GB"/></w:rPr></w:pPr>' +
    '<w:r w:rsidRPr="00107C01"><w:rPr><w:lang w:val="en-GB"/></w:rPr><w:lastRenderedPageBreak/>' +
    '<w:t>two locations adjacent to each other, first location is start location and contains all commans: look, use, use_with, take</w:t></w:r></w:p>'
Set-ParagraphOpenXml $idxTwoLocations $twoLocationsXml

$idxCreateMapping = Find-ParagraphIndex("Create mapping annotations for Jackson")
$createMappingXml = '<w:p w:rsidR="008B0220" w:rsidRPr="00107C01" w:rsidRDefault="008B0220" w:rsidP="00107C01">' +
    '<w:pPr><w:pStyle w:val="Listenabsatz"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="21"/></w:numPr>' +
    '<w:ind w:left="1068"/><w:rPr><w:lang w:val="en-GB"/></w:rPr></w:pPr>' +
    '<w:r w:rsidRPr="00107C01"><w:rPr><w:lang w:val="en-GB"/></w:rPr>' +
    '<w:t>Create mapping annotations for Jackson in all game classes, to be able to load them from JSON.</w:t></w:r></w:p>'
Set-ParagraphOpenXml $idxCreateMapping $createMappingXml

Write-Output "done"
